$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All fitness values in column C (rows 2 through 252) are corrected to a
# uniform value of 7569 as part of the SA algorithm correction.
$ws.Range("C2:C252").Value = 7569
